# Update column F ("dSF") values for the maeda_kenta workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = 7
    4  = 1
    5  = -1
    6  = 4
    7  = -1
    8  = 0
    10 = 2
    11 = -3
    12 = 3
    13 = 2
    14 = -1
    15 = 2
    17 = 2
    18 = -1
    19 = -3
    20 = -1
    21 = 1
    22 = 2
    23 = 1
    25 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
